$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 8.518499999999994
    "A8"  = -22.38590000000002
    "A10" = -21.92029999999998
    "A12" = -21.47609999999999
    "B12" = 4.835099999999997
    "B15" = 5.658899999999996
    "B17" = 4.474
    "A18" = -22.23410000000001
    "B26" = 4.299200000000002
    "B27" = 6.4714
    "B28" = 6.045399999999998
    "A37" = -21.87289999999999
    "B37" = 6.0928
    "B47" = 6.428100000000001
    "A55" = -22.0493
    "B65" = 5.935299999999995
    "A68" = -21.4333
    "B73" = 9.183600000000004
    "A77" = -20.06859999999999
    "A78" = -19.73979999999998
    "A81" = -22.22710000000001
    "A82" = -21.83830000000001
    "B84" = 5.532599999999998
    "B85" = 5.735399999999999
    "B93" = 5.434299999999997
    "B95" = 5.732800000000002
    "B98" = 5.150300000000007
    "B99" = 5.785599999999997
    "B101" = 6.116399999999998
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
